# Auto-generated edit script: updates currentAveragePrice / Leve price & profit
# columns (H-N) on specific rows across multiple worksheets, per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 1060.7333   # H28
$ws.Cells.Item(28, 9).Value = 382   # I28
$ws.Cells.Item(28, 11).Value = 382   # K28
$ws.Cells.Item(28, 13).Value = 103   # M28

# Row 107
$ws.Cells.Item(107, 8).Value = 875.125   # H107
$ws.Cells.Item(107, 9).Value = 850.1667   # I107
$ws.Cells.Item(107, 10).Value = 950   # J107
$ws.Cells.Item(107, 11).Value = 850.1667   # K107
$ws.Cells.Item(107, 12).Value = 950   # L107
$ws.Cells.Item(107, 13).Value = 1069.8333   # M107
$ws.Cells.Item(107, 14).Value = -4790   # N107

# Row 117
$ws.Cells.Item(117, 8).Value = 0   # H117
$ws.Cells.Item(117, 10).Value = 0   # J117
$ws.Cells.Item(117, 12).Value = 0   # L117
$ws.Cells.Item(117, 14).ClearContents()   # N117

$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Cells.Item(10, 8).Value = 0   # H10
$ws.Cells.Item(10, 9).Value = 0   # I10
$ws.Cells.Item(10, 11).Value = 0   # K10
$ws.Cells.Item(10, 13).ClearContents()   # M10

# Row 32
$ws.Cells.Item(32, 8).Value = 4541.1123   # H32
$ws.Cells.Item(32, 9).Value = 4109.6875   # I32
$ws.Cells.Item(32, 11).Value = 4109.6875   # K32
$ws.Cells.Item(32, 13).Value = -3822.6875   # M32

# Row 61
$ws.Cells.Item(61, 8).Value = 3888.5789   # H61
$ws.Cells.Item(61, 9).Value = 3271.0857   # I61
$ws.Cells.Item(61, 10).Value = 11092.667   # J61
$ws.Cells.Item(61, 11).Value = 3271.0857   # K61
$ws.Cells.Item(61, 12).Value = 11092.667   # L61
$ws.Cells.Item(61, 13).Value = -3059.0857   # M61
$ws.Cells.Item(61, 14).Value = -11516.667   # N61

# Row 63
$ws.Cells.Item(63, 8).Value = 7077.9473   # H63
$ws.Cells.Item(63, 10).Value = 8981   # J63
$ws.Cells.Item(63, 12).Value = 8981   # L63
$ws.Cells.Item(63, 14).Value = -10353   # N63

# Row 66
$ws.Cells.Item(66, 8).Value = 7077.9473   # H66
$ws.Cells.Item(66, 10).Value = 8981   # J66
$ws.Cells.Item(66, 12).Value = 44905   # L66
$ws.Cells.Item(66, 14).Value = -51769   # N66

# Row 88
$ws.Cells.Item(88, 8).Value = 5780.5   # H88
$ws.Cells.Item(88, 10).Value = 5230   # J88
$ws.Cells.Item(88, 12).Value = 5230   # L88
$ws.Cells.Item(88, 14).Value = -6042   # N88

# Row 91
$ws.Cells.Item(91, 8).Value = 5780.5   # H91
$ws.Cells.Item(91, 10).Value = 5230   # J91
$ws.Cells.Item(91, 12).Value = 5230   # L91
$ws.Cells.Item(91, 14).Value = -8038   # N91

# Row 92
$ws.Cells.Item(92, 8).Value = 20000   # H92
$ws.Cells.Item(92, 10).Value = 20000   # J92
$ws.Cells.Item(92, 12).Value = 20000   # L92
$ws.Cells.Item(92, 14).Value = -24992   # N92

# Row 132
$ws.Cells.Item(132, 8).Value = 1601.9445   # H132
$ws.Cells.Item(132, 9).Value = 1104.4572   # I132
$ws.Cells.Item(132, 10).Value = 19014   # J132
$ws.Cells.Item(132, 11).Value = 3313.3716   # K132
$ws.Cells.Item(132, 12).Value = 57042   # L132
$ws.Cells.Item(132, 13).Value = -783.3716000000004   # M132
$ws.Cells.Item(132, 14).Value = -62102   # N132

# Row 136
$ws.Cells.Item(136, 8).Value = 3888.5789   # H136
$ws.Cells.Item(136, 9).Value = 3271.0857   # I136
$ws.Cells.Item(136, 10).Value = 11092.667   # J136
$ws.Cells.Item(136, 11).Value = 9813.257100000001   # K136
$ws.Cells.Item(136, 12).Value = 33278.001   # L136
$ws.Cells.Item(136, 13).Value = -7263.257100000001   # M136
$ws.Cells.Item(136, 14).Value = -38378.001   # N136

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 1547.6957   # H99
$ws.Cells.Item(99, 9).Value = 1456.2106   # I99
$ws.Cells.Item(99, 11).Value = 1456.2106   # K99
$ws.Cells.Item(99, 13).Value = 41.78939999999989   # M99

# Row 134
$ws.Cells.Item(134, 8).Value = 1443.6444   # H134
$ws.Cells.Item(134, 9).Value = 1443.6444   # I134
$ws.Cells.Item(134, 10).Value = 0   # J134
$ws.Cells.Item(134, 11).Value = 4330.933199999999   # K134
$ws.Cells.Item(134, 12).Value = 0   # L134
$ws.Cells.Item(134, 13).Value = -1795.933199999999   # M134
$ws.Cells.Item(134, 14).ClearContents()   # N134

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 33772.5   # H31
$ws.Cells.Item(31, 9).Value = 1325.85   # I31
$ws.Cells.Item(31, 11).Value = 1325.85   # K31
$ws.Cells.Item(31, 13).Value = -1030.85   # M31

# Row 34
$ws.Cells.Item(34, 8).Value = 33772.5   # H34
$ws.Cells.Item(34, 9).Value = 1325.85   # I34
$ws.Cells.Item(34, 11).Value = 1325.85   # K34
$ws.Cells.Item(34, 13).Value = -1123.85   # M34

# Row 58
$ws.Cells.Item(58, 8).Value = 3598.258   # H58
$ws.Cells.Item(58, 9).Value = 1583.3334   # I58
$ws.Cells.Item(58, 10).Value = 7829.6   # J58
$ws.Cells.Item(58, 11).Value = 1583.3334   # K58
$ws.Cells.Item(58, 12).Value = 7829.6   # L58
$ws.Cells.Item(58, 13).Value = -1380.3334   # M58
$ws.Cells.Item(58, 14).Value = -8235.6   # N58

# Row 107
$ws.Cells.Item(107, 8).Value = 1975.7   # H107
$ws.Cells.Item(107, 9).Value = 1309.7858   # I107
$ws.Cells.Item(107, 10).Value = 3529.5   # J107
$ws.Cells.Item(107, 11).Value = 1309.7858   # K107
$ws.Cells.Item(107, 12).Value = 3529.5   # L107
$ws.Cells.Item(107, 13).Value = 610.2141999999999   # M107
$ws.Cells.Item(107, 14).Value = -7369.5   # N107

# Row 132
$ws.Cells.Item(132, 8).Value = 3041   # H132
$ws.Cells.Item(132, 9).Value = 2364.8572   # I132
$ws.Cells.Item(132, 11).Value = 7094.571599999999   # K132
$ws.Cells.Item(132, 13).Value = -4564.571599999999   # M132

# Row 134
$ws.Cells.Item(134, 8).Value = 1658.0769   # H134
$ws.Cells.Item(134, 9).Value = 1306.6052   # I134
$ws.Cells.Item(134, 11).Value = 3919.8156   # K134
$ws.Cells.Item(134, 13).Value = -1384.8156   # M134

# Row 136
$ws.Cells.Item(136, 8).Value = 3598.258   # H136
$ws.Cells.Item(136, 9).Value = 1583.3334   # I136
$ws.Cells.Item(136, 10).Value = 7829.6   # J136
$ws.Cells.Item(136, 11).Value = 4750.0002   # K136
$ws.Cells.Item(136, 12).Value = 23488.8   # L136
$ws.Cells.Item(136, 13).Value = -2200.0002   # M136
$ws.Cells.Item(136, 14).Value = -28588.8   # N136

$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Cells.Item(109, 8).Value = 2256.1667   # H109
$ws.Cells.Item(109, 9).Value = 1165.75   # I109
$ws.Cells.Item(109, 11).Value = 3497.25   # K109
$ws.Cells.Item(109, 13).Value = -2457.25   # M109

# Row 132
$ws.Cells.Item(132, 8).Value = 4580.5806   # H132
$ws.Cells.Item(132, 9).Value = 3412.353   # I132
$ws.Cells.Item(132, 11).Value = 30711.177   # K132
$ws.Cells.Item(132, 13).Value = -28181.177   # M132

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 512.5   # H2
$ws.Cells.Item(2, 10).Value = 1358.4286   # J2
$ws.Cells.Item(2, 12).Value = 1358.4286   # L2
$ws.Cells.Item(2, 14).Value = -1584.4286   # N2

# Row 80
$ws.Cells.Item(80, 8).Value = 841002   # H80
$ws.Cells.Item(80, 9).Value = 5000000   # I80
$ws.Cells.Item(80, 10).Value = 9202.4   # J80
$ws.Cells.Item(80, 11).Value = 5000000   # K80
$ws.Cells.Item(80, 12).Value = 9202.4   # L80
$ws.Cells.Item(80, 13).Value = -4999002   # M80
$ws.Cells.Item(80, 14).Value = -11198.4   # N80

# Row 83
$ws.Cells.Item(83, 8).Value = 841002   # H83
$ws.Cells.Item(83, 9).Value = 5000000   # I83
$ws.Cells.Item(83, 10).Value = 9202.4   # J83
$ws.Cells.Item(83, 11).Value = 25000000   # K83
$ws.Cells.Item(83, 12).Value = 46012   # L83
$ws.Cells.Item(83, 13).Value = -24995008   # M83
$ws.Cells.Item(83, 14).Value = -55996   # N83

# Row 103
$ws.Cells.Item(103, 8).Value = 31051.666   # H103
$ws.Cells.Item(103, 10).Value = 31051.666   # J103
$ws.Cells.Item(103, 12).Value = 31051.666   # L103
$ws.Cells.Item(103, 14).Value = -33395.666   # N103

# Row 113
$ws.Cells.Item(113, 8).Value = 4069.7144   # H113
$ws.Cells.Item(113, 9).Value = 3426.5715   # I113
$ws.Cells.Item(113, 10).Value = 4712.857   # J113
$ws.Cells.Item(113, 11).Value = 3426.5715   # K113
$ws.Cells.Item(113, 12).Value = 4712.857   # L113
$ws.Cells.Item(113, 13).Value = -1256.5715   # M113
$ws.Cells.Item(113, 14).Value = -9052.857   # N113

# Row 126
$ws.Cells.Item(126, 8).Value = 3818.16   # H126
$ws.Cells.Item(126, 9).Value = 2427   # I126
$ws.Cells.Item(126, 10).Value = 4911.2144   # J126
$ws.Cells.Item(126, 11).Value = 7281   # K126
$ws.Cells.Item(126, 12).Value = 14733.6432   # L126
$ws.Cells.Item(126, 13).Value = -4811   # M126
$ws.Cells.Item(126, 14).Value = -19673.6432   # N126

# Row 132
$ws.Cells.Item(132, 8).Value = 2115.7314   # H132
$ws.Cells.Item(132, 9).Value = 1824.2222   # I132
$ws.Cells.Item(132, 11).Value = 5472.6666   # K132
$ws.Cells.Item(132, 13).Value = -2942.6666   # M132

# Row 136
$ws.Cells.Item(136, 8).Value = 63498.363   # H136
$ws.Cells.Item(136, 10).Value = 63498.363   # J136
$ws.Cells.Item(136, 12).Value = 190495.089   # L136
$ws.Cells.Item(136, 14).Value = -195595.089   # N136

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 13888.5   # H7
$ws.Cells.Item(7, 10).Value = 23833.334   # J7
$ws.Cells.Item(7, 12).Value = 23833.334   # L7
$ws.Cells.Item(7, 14).Value = -24057.334   # N7

# Row 12
$ws.Cells.Item(12, 8).Value = 10000   # H12
$ws.Cells.Item(12, 9).Value = 10000   # I12
$ws.Cells.Item(12, 10).Value = 0   # J12
$ws.Cells.Item(12, 11).Value = 10000   # K12
$ws.Cells.Item(12, 12).Value = 0   # L12
$ws.Cells.Item(12, 13).Value = -9830   # M12
$ws.Cells.Item(12, 14).ClearContents()   # N12

# Row 40
$ws.Cells.Item(40, 8).Value = 13128   # H40
$ws.Cells.Item(40, 9).Value = 22643.334   # I40
$ws.Cells.Item(40, 11).Value = 22643.334   # K40
$ws.Cells.Item(40, 13).Value = -22507.334   # M40

# Row 93
$ws.Cells.Item(93, 8).Value = 2355.6538   # H93
$ws.Cells.Item(93, 9).Value = 2351.5   # I93
$ws.Cells.Item(93, 11).Value = 2351.5   # K93
$ws.Cells.Item(93, 13).Value = -1103.5   # M93

# Row 126
$ws.Cells.Item(126, 8).Value = 13888.5   # H126
$ws.Cells.Item(126, 10).Value = 23833.334   # J126
$ws.Cells.Item(126, 12).Value = 71500.00199999999   # L126
$ws.Cells.Item(126, 14).Value = -76440.00199999999   # N126

# Row 132
$ws.Cells.Item(132, 8).Value = 9459.083000000001   # H132
$ws.Cells.Item(132, 9).Value = 8000.1665   # I132
$ws.Cells.Item(132, 10).Value = 10918   # J132
$ws.Cells.Item(132, 11).Value = 24000.4995   # K132
$ws.Cells.Item(132, 12).Value = 32754   # L132
$ws.Cells.Item(132, 13).Value = -21470.4995   # M132
$ws.Cells.Item(132, 14).Value = -37814   # N132

# Row 136
$ws.Cells.Item(136, 8).Value = 6009.7896   # H136
$ws.Cells.Item(136, 9).Value = 4876.231   # I136
$ws.Cells.Item(136, 10).Value = 8465.833000000001   # J136
$ws.Cells.Item(136, 11).Value = 14628.693   # K136
$ws.Cells.Item(136, 12).Value = 25397.499   # L136
$ws.Cells.Item(136, 13).Value = -12078.693   # M136
$ws.Cells.Item(136, 14).Value = -30497.499   # N136

$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Cells.Item(12, 8).Value = 0   # H12
$ws.Cells.Item(12, 10).Value = 0   # J12
$ws.Cells.Item(12, 12).Value = 0   # L12
$ws.Cells.Item(12, 14).ClearContents()   # N12

# Row 113
$ws.Cells.Item(113, 8).Value = 379.875   # H113
$ws.Cells.Item(113, 9).Value = 357.57895   # I113
$ws.Cells.Item(113, 11).Value = 1072.73685   # K113
$ws.Cells.Item(113, 13).Value = 1097.26315   # M113

# Row 122
$ws.Cells.Item(122, 8).Value = 2832.2856   # H122
$ws.Cells.Item(122, 9).Value = 2095.3809   # I122
$ws.Cells.Item(122, 11).Value = 6286.1427   # K122
$ws.Cells.Item(122, 13).Value = -3836.1427   # M122

# Row 126
$ws.Cells.Item(126, 8).Value = 3636.9583   # H126
$ws.Cells.Item(126, 9).Value = 2541.4211   # I126
$ws.Cells.Item(126, 11).Value = 7624.263300000001   # K126
$ws.Cells.Item(126, 13).Value = -5154.263300000001   # M126

# Row 132
$ws.Cells.Item(132, 8).Value = 1824.8379   # H132
$ws.Cells.Item(132, 9).Value = 1347.6111   # I132
$ws.Cells.Item(132, 11).Value = 4042.8333   # K132
$ws.Cells.Item(132, 13).Value = -1512.8333   # M132

# Row 136
$ws.Cells.Item(136, 8).Value = 2886.25   # H136
$ws.Cells.Item(136, 9).Value = 1213.5217   # I136
$ws.Cells.Item(136, 10).Value = 10580.8   # J136
$ws.Cells.Item(136, 11).Value = 3640.5651   # K136
$ws.Cells.Item(136, 12).Value = 31742.4   # L136
$ws.Cells.Item(136, 13).Value = -1090.5651   # M136
$ws.Cells.Item(136, 14).Value = -36842.39999999999   # N136

